$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '34.247.07'
$ws.Range("E2").Value = '  +0.46%  '

# Row 3
$ws.Range("D3").Value = '1.791.06'
$ws.Range("E3").Value = '  +0.14%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.22'
$ws.Range("E5").Value = '  -0.23%  '

# Row 6
$ws.Range("E6").Value = '  +0.36%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.27'
$ws.Range("E8").Value = '  +0.49%  '

# Row 9
$ws.Range("E9").Value = '  +0.29%  '

# Row 10
$ws.Range("E10").Value = '  -0.40%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0946'
$ws.Range("E11").Value = '  +0.71%  '

# Row 12
$ws.Range("D12").Value = '2.047.13'
$ws.Range("E12").Value = '  +0.02%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.14'
$ws.Range("E13").Value = '  -3.59%  '

# Row 14
$ws.Range("D14").Value = '1.776.05'
$ws.Range("E14").Value = '  -0.80%  '

# Row 15
$ws.Range("E15").Value = '  +0.62%  '

# Row 16
$ws.Range("D16").Value = '34.228.95'
$ws.Range("E16").Value = '  +0.44%  '

# Row 17
$ws.Range("E17").Value = '  +0.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.06'
$ws.Range("E18").Value = '  +0.09%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0808'
$ws.Range("E19").Value = '  +3.30%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '246.72'
$ws.Range("E20").Value = '  +0.86%  '

# Row 21
$ws.Range("E21").Value = '  +0.58%  '

# Row 22
$ws.Range("E22").Value = '  +0.02%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.21'
$ws.Range("E23").Value = '  +2.60%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  -0.15%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.55'
$ws.Range("E25").Value = '  -0.63%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.19'
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.35'
$ws.Range("E27").Value = '  +0.37%  '

# Row 28
$ws.Range("E28").Value = '  +0.91%  '

# Row 29
$ws.Range("E29").Value = '  +0.15%  '

# Row 30
$ws.Range("E30").Value = '  +0.03%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0521'
$ws.Range("E31").Value = '  +0.05%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.76'
$ws.Range("E32").Value = '  +2.66%  '

# Row 33
$ws.Range("E33").Value = '  +4.09%  '

# Row 34
$ws.Range("E34").Value = '  -1.34%  '

# Row 35
$ws.Range("D35").Value = '1.444.03'
$ws.Range("E35").Value = '  +1.98%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.58'
$ws.Range("E36").Value = '  +9.79%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.670'
$ws.Range("E37").Value = '  +4.09%  '

# Row 38
$ws.Range("E38").Value = '  +0.04%  '

# Row 39
$ws.Range("E39").Value = '  +1.32%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '82.09'
$ws.Range("E40").Value = '  +1.93%  '

# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.11'
$ws.Range("E41").Value = '  +5.32%  '

# Row 42
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.38'
$ws.Range("E42").Value = '  +1.13%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.923'
$ws.Range("E43").Value = '  +0.61%  '

# Row 44
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.72'
$ws.Range("E44").Value = '  +1.33%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0519'
$ws.Range("E45").Value = '  +2.20%  '

# Row 46
$ws.Range("E46").Value = '  +0.41%  '

# Row 47
$ws.Range("E47").Value = '  +0.22%  '

# Row 48
$ws.Range("D48").Value = '1.946.90'
$ws.Range("E48").Value = '  -0.06%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.53'
$ws.Range("E49").Value = '  -1.56%  '

# Row 50
$ws.Range("E50").Value = '  -0.04%  '

# Row 51
$ws.Range("D51").Value = '0.0₆0126'
$ws.Range("E51").Value = '  -7.87%  '
